$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting existing C:L to D:M
$ws.Columns("C").Insert()

# New header for the inserted column
$ws.Range("C1").Value = "statut_name"

# Fill the new column with the constant status text for each data row (rows 2-11)
$statutName = "pas de résultat ni de publication"
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = $statutName
}
